# Update the "想去人数" (number of people interested) column (F) for several
# rows in both the "展览" and "全部类型" worksheets, which carry duplicated data.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 627
    5  = 29
    13 = 1101
    25 = 1670
    31 = 3928
    34 = 226
    35 = 994
    40 = 45
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
